$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.228
$ws.Range("A9").Value = -20.912
$ws.Range("B11").Value = 6.927
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
$ws.Range("E21").Value = 13.186
